$d = $word.ActiveDocument

# This template contains two "field code" style paragraphs built from
# w:fldChar / w:instrText runs (Word field codes). The new parser expects
# plain M2Doc template tags written directly as literal text
# ("{m:...}") instead of Word fields, so we rewrite each of those two
# paragraphs, replacing the fldChar/instrText runs with plain w:t runs
# that spell out the same text wrapped in curly braces.

$pkgHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 2: {m:'Some value'.setDocumentCompany()} ---
$p2 = $d.Paragraphs.Item(2)
$p2Xml = '<w:p w14:paraId="5F0A223D" w14:textId="45C757E4" w:rsidR="00E1471F" w:rsidRDefault="00E1471F" w:rsidP="00E1471F">' +
  '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:t>{m:</w:t></w:r>' +
  '<w:r><w:t>' + [char]39 + '</w:t></w:r>' +
  '<w:r><w:t>Some value</w:t></w:r>' +
  '<w:r><w:t>' + [char]39 + '</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '<w:r><w:t>setDocument</w:t></w:r>' +
  '<w:r><w:t>Company</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">()}</w:t></w:r>' +
  '</w:p>'
$p2.Range.InsertXML($pkgHeader + $p2Xml + $pkgFooter) | Out-Null

# --- Paragraph 3: {m:''.getDocumentCompany()} ---
$p3 = $d.Paragraphs.Item(3)
$p3Xml = '<w:p w14:paraId="2C980985" w14:textId="1C6286AE" w:rsidR="00CD75A1" w:rsidRDefault="00CD75A1" w:rsidP="00CD75A1">' +
  '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:t>{m:' + [char]39 + [char]39 + '.g</w:t></w:r>' +
  '<w:r><w:t>etDocument</w:t></w:r>' +
  '<w:r><w:t>C</w:t></w:r>' +
  '<w:r><w:t>ompany</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">()}</w:t></w:r>' +
  '</w:p>'
$p3.Range.InsertXML($pkgHeader + $p3Xml + $pkgFooter) | Out-Null
